# Updated codes in Financial page. Added new test cases.
# This script reproduces, via Excel COM automation, the data/formatting
# changes described in the target diff for xl/worksheets/sheet1.xml,
# xl/sharedStrings.xml and xl/styles.xml of the "admin" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlLeft = -4131

# ------------------------------------------------------------------
# Build a helper cell far off-sheet that carries the *new* cell style
# (existing bordered style plus left horizontal alignment). Copying an
# existing bordered cell (A2, style index 3) first guarantees the new
# style reuses the same border definition instead of Excel creating a
# brand-new (duplicate) border entry.
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)
$ws.Range("Z1").HorizontalAlignment = $xlLeft

# ------------------------------------------------------------------
# Row 4 - replace the previous "test_1 / order / order value" content
# with the new order-related test row, now using the bordered style
# (same as rows 2-3) instead of the default style.
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A4:C4").PasteSpecial($xlPasteFormats)
$ws.Range("A4").Value = "test_order"
$ws.Range("B4").Value = "order"
$ws.Range("C4").Value = "Server"

# ------------------------------------------------------------------
# Row 5 - test_selectOrderType / projectid / 98313 (numeric, left align)
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A5:B5").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Value = "test_selectOrderType"
$ws.Range("B5").Value = "projectid"
$ws.Range("Z1").Copy()
$ws.Range("C5").PasteSpecial($xlPasteFormats)
$ws.Range("C5").Value = 98313

# ------------------------------------------------------------------
# Row 6 - test_projectinfo / country / US - 2020
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A6:C6").PasteSpecial($xlPasteFormats)
$ws.Range("A6").Value = "test_projectinfo"
$ws.Range("B6").Value = "country"
$ws.Range("C6").Value = "US - 2020"

# ------------------------------------------------------------------
# Row 7 - test_projectinfo / Capital Amount / 123 (numeric, left align)
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A7:B7").PasteSpecial($xlPasteFormats)
$ws.Range("A7").Value = "test_projectinfo"
$ws.Range("B7").Value = "Capital Amount"
$ws.Range("Z1").Copy()
$ws.Range("C7").PasteSpecial($xlPasteFormats)
$ws.Range("C7").Value = 123

# ------------------------------------------------------------------
# Row 8 - test_projectinfo / Expense Amount / 1 (numeric, left align)
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A8:B8").PasteSpecial($xlPasteFormats)
$ws.Range("A8").Value = "test_projectinfo"
$ws.Range("B8").Value = "Expense Amount"
$ws.Range("Z1").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)
$ws.Range("C8").Value = 1

# ------------------------------------------------------------------
# Row 9 - test_projectinfo / Quote Number / test123
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A9:B9").PasteSpecial($xlPasteFormats)
$ws.Range("A9").Value = "test_projectinfo"
$ws.Range("B9").Value = "Quote Number"
$ws.Range("Z1").Copy()
$ws.Range("C9").PasteSpecial($xlPasteFormats)
$ws.Range("C9").Value = "test123"

# ------------------------------------------------------------------
# Row 10 - test_projectinfo / Supplier / 3D NETWORKS
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A10:B10").PasteSpecial($xlPasteFormats)
$ws.Range("A10").Value = "test_projectinfo"
$ws.Range("B10").Value = "Supplier"
$ws.Range("Z1").Copy()
$ws.Range("C10").PasteSpecial($xlPasteFormats)
$ws.Range("C10").Value = "3D NETWORKS"

# ------------------------------------------------------------------
# Row 11 - test_projectinfo / upload / Dummyfor3PI
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A11:B11").PasteSpecial($xlPasteFormats)
$ws.Range("A11").Value = "test_projectinfo"
$ws.Range("B11").Value = "upload"
$ws.Range("Z1").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)
$ws.Range("C11").Value = "Dummyfor3PI"

# ------------------------------------------------------------------
# Row 12 - test_projectinfo / path / C:\Users\satheeshnair\Desktop\word\Dummyfor3PI.txt
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A12:B12").PasteSpecial($xlPasteFormats)
$ws.Range("A12").Value = "test_projectinfo"
$ws.Range("B12").Value = "path"
$ws.Range("Z1").Copy()
$ws.Range("C12").PasteSpecial($xlPasteFormats)
$ws.Range("C12").Value = "C:\Users\satheeshnair\Desktop\word\Dummyfor3PI.txt"

# ------------------------------------------------------------------
# Rows 13-14 - blank spacer rows that keep the bordered styling
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A13:B14").PasteSpecial($xlPasteFormats)
$ws.Range("Z1").Copy()
$ws.Range("C13:C14").PasteSpecial($xlPasteFormats)

# Remove the temporary helper cell so it does not show up in the sheet.
$ws.Range("Z1").Clear()

# ------------------------------------------------------------------
# Column C is now much wider to fit the long file-path text.
# ------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 55.3

# ------------------------------------------------------------------
# Update the active selection to match the new cursor position.
# ------------------------------------------------------------------
$ws.Range("F11").Select()
